$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 96
$ws.Range("H96").Value = 287.77777
$ws.Range("J96").Value = 399
$ws.Range("L96").Value = 1197
$ws.Range("N96").Value = -3943
# Row 100
$ws.Range("H100").Value = 5777.5293
$ws.Range("I100").Value = 2672.8572
$ws.Range("J100").Value = 7950.8
$ws.Range("K100").Value = 2672.8572
$ws.Range("L100").Value = 7950.8
$ws.Range("M100").Value = -2131.8572
$ws.Range("N100").Value = -9032.799999999999
# Row 137
$ws.Range("H137").Value = 2530.5789
$ws.Range("I137").Value = 1788.9333
$ws.Range("K137").Value = 5366.7999
$ws.Range("M137").Value = -2816.7999
# Row 141
$ws.Range("H141").Value = 7070
$ws.Range("I141").Value = 7508
$ws.Range("K141").Value = 22524
$ws.Range("M141").Value = -17344

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 10564.254
$ws.Range("I32").Value = 8163.778
$ws.Range("K32").Value = 8163.778
$ws.Range("M32").Value = -7876.778
# Row 88
$ws.Range("H88").Value = 2004.7826
$ws.Range("I88").Value = 2005.3334
$ws.Range("J88").Value = 2004.1818
$ws.Range("K88").Value = 2005.3334
$ws.Range("L88").Value = 2004.1818
$ws.Range("M88").Value = -1599.3334
$ws.Range("N88").Value = -2816.1818
# Row 91
$ws.Range("H91").Value = 2004.7826
$ws.Range("I91").Value = 2005.3334
$ws.Range("J91").Value = 2004.1818
$ws.Range("K91").Value = 2005.3334
$ws.Range("L91").Value = 2004.1818
$ws.Range("M91").Value = -601.3334
$ws.Range("N91").Value = -4812.1818
# Row 110
$ws.Range("H110").Value = 1519.0322
$ws.Range("I110").Value = 1668.9615
$ws.Range("K110").Value = 1668.9615
$ws.Range("M110").Value = 376.0385000000001
# Row 119
$ws.Range("H119").Value = 89998.75
$ws.Range("J119").Value = 89998.75
$ws.Range("L119").Value = 89998.75
$ws.Range("N119").Value = -99674.75
# Row 128
$ws.Range("H128").Value = 62057.4
$ws.Range("J128").Value = 62057.4
$ws.Range("L128").Value = 62057.4
$ws.Range("N128").Value = -72017.39999999999
# Row 131
$ws.Range("H131").Value = 68452.5
$ws.Range("J131").Value = 68452.5
$ws.Range("L131").Value = 68452.5
$ws.Range("N131").Value = -78532.5
# Row 132
$ws.Range("H132").Value = 52636680
$ws.Range("I132").Value = 90913540
$ws.Range("K132").Value = 272740620
$ws.Range("M132").Value = -272738090

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 5444.5
$ws.Range("I86").Value = 1372.375
$ws.Range("J86").Value = 13588.75
$ws.Range("K86").Value = 1372.375
$ws.Range("L86").Value = 13588.75
$ws.Range("M86").Value = -249.375
$ws.Range("N86").Value = -15834.75
# Row 89
$ws.Range("H89").Value = 5444.5
$ws.Range("I89").Value = 1372.375
$ws.Range("J89").Value = 13588.75
$ws.Range("K89").Value = 6861.875
$ws.Range("L89").Value = 67943.75
$ws.Range("M89").Value = -1245.875
$ws.Range("N89").Value = -79175.75
# Row 134
$ws.Range("H134").Value = 5890.8335
$ws.Range("I134").Value = 4540.125
$ws.Range("K134").Value = 13620.375
$ws.Range("M134").Value = -11085.375

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4787.8623
$ws.Range("I31").Value = 1394.2
$ws.Range("K31").Value = 1394.2
$ws.Range("M31").Value = -1099.2
# Row 34
$ws.Range("H34").Value = 4787.8623
$ws.Range("I34").Value = 1394.2
$ws.Range("K34").Value = 1394.2
$ws.Range("M34").Value = -1192.2
# Row 63
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
# Row 66
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
# Row 122
$ws.Range("H122").Value = 43481070
$ws.Range("I122").Value = 76924380
$ws.Range("K122").Value = 230773140
$ws.Range("M122").Value = -230770690
# Row 134
$ws.Range("H134").Value = 5327.55
$ws.Range("I134").Value = 4305.625
$ws.Range("J134").Value = 6860.4375
$ws.Range("K134").Value = 12916.875
$ws.Range("L134").Value = 20581.3125
$ws.Range("M134").Value = -10381.875
$ws.Range("N134").Value = -25651.3125

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 97
$ws.Range("H97").Value = 1019.04
$ws.Range("J97").Value = 480.75
$ws.Range("L97").Value = 1442.25
$ws.Range("N97").Value = -2434.25
# Row 131
$ws.Range("H131").Value = 6240.7334
$ws.Range("J131").Value = 8192
$ws.Range("L131").Value = 24576
$ws.Range("N131").Value = -34656

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 38098292
$ws.Range("I102").Value = 5496565.5
$ws.Range("J102").Value = 250009500
$ws.Range("K102").Value = 5496565.5
$ws.Range("L102").Value = 250009500
$ws.Range("M102").Value = -5494943.5
$ws.Range("N102").Value = -250012744
# Row 105
$ws.Range("H105").Value = 66492.16
$ws.Range("J105").Value = 66492.16
$ws.Range("L105").Value = 66492.16
$ws.Range("N105").Value = -73480.16
# Row 113
$ws.Range("H113").Value = 11640.728
$ws.Range("I113").Value = 6609.8
$ws.Range("K113").Value = 6609.8
$ws.Range("M113").Value = -4439.8
# Row 132
$ws.Range("H132").Value = 4844.75
$ws.Range("I132").Value = 3888.8823
$ws.Range("J132").Value = 6322
$ws.Range("K132").Value = 11666.6469
$ws.Range("L132").Value = 18966
$ws.Range("M132").Value = -9136.6469
$ws.Range("N132").Value = -24026

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 1132.3334
$ws.Range("I55").Value = 1532.6666
$ws.Range("K55").Value = 1532.6666
$ws.Range("M55").Value = -1359.6666
# Row 61
$ws.Range("H61").Value = 2226.4443
$ws.Range("I61").Value = 1573.0714
$ws.Range("K61").Value = 1573.0714
$ws.Range("M61").Value = -1371.0714
# Row 76
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
# Row 79
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
# Row 93
$ws.Range("H93").Value = 387961.56
$ws.Range("I93").Value = 3105.4211
$ws.Range("K93").Value = 3105.4211
$ws.Range("M93").Value = -1857.4211
# Row 97
$ws.Range("H97").Value = 39672
$ws.Range("J97").Value = 39672
$ws.Range("L97").Value = 39672
$ws.Range("N97").Value = -41654
# Row 113
$ws.Range("H113").Value = 2226.4443
$ws.Range("I113").Value = 1573.0714
$ws.Range("K113").Value = 1573.0714
$ws.Range("M113").Value = 596.9286

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 304.9
$ws.Range("I107").Value = 294.85715
$ws.Range("K107").Value = 884.5714499999999
$ws.Range("M107").Value = 1035.42855
# Row 122
$ws.Range("H122").Value = 18446.166
$ws.Range("I122").Value = 21001.732
$ws.Range("J122").Value = 5668.3335
$ws.Range("K122").Value = 63005.196
$ws.Range("L122").Value = 17005.0005
$ws.Range("M122").Value = -60555.196
$ws.Range("N122").Value = -21905.0005
# Row 130
$ws.Range("H130").Value = 50429
$ws.Range("J130").Value = 50429
$ws.Range("L130").Value = 50429
$ws.Range("N130").Value = -60469
